# Update "想去人数" (interest count) figures across the workbook sheets.
# This mirrors a re-scrape of the source data where the F column counters
# changed slightly for a number of rows on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 269
$ws1.Range("F8").Value  = 330
$ws1.Range("F9").Value  = 2191
$ws1.Range("F10").Value = 1155
$ws1.Range("F11").Value = 1043
$ws1.Range("F12").Value = 844
$ws1.Range("F14").Value = 837
$ws1.Range("F15").Value = 1464
$ws1.Range("F16").Value = 702
$ws1.Range("F17").Value = 1686
$ws1.Range("F18").Value = 39
$ws1.Range("F19").Value = 347
$ws1.Range("F20").Value = 63
$ws1.Range("F21").Value = 98

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 24
$ws2.Range("F19").Value = 151
$ws2.Range("F24").Value = 87
$ws2.Range("F28").Value = 180
$ws2.Range("F38").Value = 337

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value  = 2480
$ws3.Range("F5").Value  = 710
$ws3.Range("F7").Value  = 9531
$ws3.Range("F8").Value  = 1842
$ws3.Range("F9").Value  = 132
$ws3.Range("F12").Value = 354
$ws3.Range("F13").Value = 2767
$ws3.Range("F14").Value = 355
$ws3.Range("F15").Value = 661

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2480
$ws4.Range("F3").Value  = 710
$ws4.Range("F4").Value  = 132
$ws4.Range("F8").Value  = 2767
$ws4.Range("F9").Value  = 355
$ws4.Range("F11").Value = 661
$ws4.Range("F17").Value = 269
$ws4.Range("F18").Value = 330
$ws4.Range("F20").Value = 1043
$ws4.Range("F21").Value = 844
$ws4.Range("F23").Value = 837
$ws4.Range("F28").Value = 702
$ws4.Range("F31").Value = 1686
$ws4.Range("F32").Value = 347
$ws4.Range("F33").Value = 87
$ws4.Range("F37").Value = 180
$ws4.Range("F39").Value = 63
$ws4.Range("F42").Value = 337
